$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix tiny floating-point adjustment on the existing last row (A15)
$ws.Range("A15").Value = 45877.58355226852

# Append the new row of sensor data (row 16)
$ws.Range("A16").Value = 45877.62522984331
$ws.Range("B16").Value = 2025
$ws.Range("C16").Value = 32
$ws.Range("D16").Value = 19.22
$ws.Range("E16").Value = 77.90000000000001
$ws.Range("F16").Value = 442.97
$ws.Range("G16").Value = 11.84
$ws.Range("H16").Value = "SE"
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = "15:00:19"

# Match formatting of the row above (style index 2 applied to column A date cells)
$ws.Range("A16").NumberFormat = $ws.Range("A15").NumberFormat
